$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V data between rows 13 and 15 (two matches reordered) ---
# Row 13
$ws.Range("F13").Value = "Lisen"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = "Jihlava"
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 1.72
$ws.Range("K13").Value = "27/07/2023 09:12"
$ws.Range("L13").Value = 1.83
$ws.Range("M13").Value = "29/07/2023 14:58"
$ws.Range("N13").Value = 3.61
$ws.Range("O13").Value = "27/07/2023 09:12"
$ws.Range("P13").Value = 3.62
$ws.Range("Q13").Value = "29/07/2023 16:51"
$ws.Range("R13").Value = 4.1
$ws.Range("S13").Value = "27/07/2023 09:12"
$ws.Range("T13").Value = 4.23
$ws.Range("U13").Value = "29/07/2023 16:51"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/czech-republic/fnl/lisen-jihlava/GQgBmAHo/"

# Row 15
$ws.Range("F15").Value = "Vlasim"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = "Chrudim"
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1.76
$ws.Range("K15").Value = "27/07/2023 09:12"
$ws.Range("L15").Value = 1.59
$ws.Range("M15").Value = "29/07/2023 16:52"
$ws.Range("N15").Value = 3.68
$ws.Range("O15").Value = "27/07/2023 09:12"
$ws.Range("P15").Value = 4.23
$ws.Range("Q15").Value = "29/07/2023 16:52"
$ws.Range("R15").Value = 3.81
$ws.Range("S15").Value = "27/07/2023 09:12"
$ws.Range("T15").Value = 5.21
$ws.Range("U15").Value = "29/07/2023 16:52"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/czech-republic/fnl/vlasim-chrudim/p6oSqSGA/"

# Row 108
$ws.Range("F108").Value = "Chrudim"
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = "Kromeriz"
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 1.93
$ws.Range("K108").Value = "26/10/2023 22:42"
$ws.Range("L108").Value = 1.87
$ws.Range("M108").Value = "28/10/2023 10:06"
$ws.Range("N108").Value = 3.41
$ws.Range("O108").Value = "26/10/2023 22:42"
$ws.Range("P108").Value = 3.48
$ws.Range("Q108").Value = "28/10/2023 10:12"
$ws.Range("R108").Value = 3.41
$ws.Range("S108").Value = "26/10/2023 22:42"
$ws.Range("T108").Value = 4.25
$ws.Range("U108").Value = "28/10/2023 10:06"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/czech-republic/fnl/chrudim-kromeriz/U7jRpNJf/"

# Row 109
$ws.Range("F109").Value = "Vyskov"
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = "Opava"
$ws.Range("I109").Value = 3
$ws.Range("J109").Value = 1.9
$ws.Range("K109").Value = "26/10/2023 22:42"
$ws.Range("L109").Value = 1.87
$ws.Range("M109").Value = "28/10/2023 10:06"
$ws.Range("N109").Value = 3.71
$ws.Range("O109").Value = "26/10/2023 22:42"
$ws.Range("P109").Value = 3.46
$ws.Range("Q109").Value = "28/10/2023 10:06"
$ws.Range("R109").Value = 3.36
$ws.Range("S109").Value = "26/10/2023 22:42"
$ws.Range("T109").Value = 4.28
$ws.Range("U109").Value = "28/10/2023 10:06"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/czech-republic/fnl/mfk-vyskov-opava/GOhwrL3D/"

# --- Append two new match rows (122, 123), copying formatting from row 121 ---
$ws.Range("A121:V121").Copy()
$ws.Range("A122:V123").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 122
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "czech-republic"
$ws.Range("C122").Value = "fnl"
$ws.Range("D122").Value = "2023-2024"
$ws.Range("E122").Value = 45240.70833333334
$ws.Range("F122").Value = "Jihlava"
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = "Lisen"
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 2.12
$ws.Range("K122").Value = "08/11/2023 17:13"
$ws.Range("L122").Value = 2.39
$ws.Range("M122").Value = "10/11/2023 16:57"
$ws.Range("N122").Value = 3.31
$ws.Range("O122").Value = "08/11/2023 17:13"
$ws.Range("P122").Value = 3.05
$ws.Range("Q122").Value = "10/11/2023 16:57"
$ws.Range("R122").Value = 3.04
$ws.Range("S122").Value = "08/11/2023 17:13"
$ws.Range("T122").Value = 3.21
$ws.Range("U122").Value = "10/11/2023 16:57"
$ws.Range("V122").Value = "https://www.betexplorer.com/football/czech-republic/fnl/jihlava-lisen/lYbSDLms/"

# Row 123
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "czech-republic"
$ws.Range("C123").Value = "fnl"
$ws.Range("D123").Value = "2023-2024"
$ws.Range("E123").Value = 45240.75
$ws.Range("F123").Value = "Pribram"
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = "Zizkov"
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 2.19
$ws.Range("K123").Value = "08/11/2023 18:12"
$ws.Range("L123").Value = 2.2
$ws.Range("M123").Value = "10/11/2023 17:59"
$ws.Range("N123").Value = 3.69
$ws.Range("O123").Value = "08/11/2023 18:12"
$ws.Range("P123").Value = 3.43
$ws.Range("Q123").Value = "10/11/2023 17:59"
$ws.Range("R123").Value = 2.75
$ws.Range("S123").Value = "08/11/2023 18:12"
$ws.Range("T123").Value = 3.21
$ws.Range("U123").Value = "10/11/2023 17:59"
$ws.Range("V123").Value = "https://www.betexplorer.com/football/czech-republic/fnl/pribram-zizkov/40NeJsQJ/"

